$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-05 Friday" "2024-07-06 Saturday"
Replace-Text "957×6=5742" "584×5=2920"
Replace-Text "555×7=3885" "299×9=2691"
Replace-Text "355×4=1420" "488×7=3416"
Replace-Text "638×6=3828" "293×5=1465"
Replace-Text "129×9=1161" "208×7=1456"
Replace-Text "399×7=2793" "365×9=3285"
Replace-Text "523×6=3138" "236×9=2124"
Replace-Text "686×4=2744" "529×4=2116"
Replace-Text "773×8=6184" "398×7=2786"
Replace-Text "657×7=4599" "605×7=4235"
Replace-Text "369×8=2952" "740×5=3700"
Replace-Text "440×5=2200" "182×4=728"
Replace-Text "204×4=816" "163×6=978"
Replace-Text "180×7=1260" "898×6=5388"
Replace-Text "752×5=3760" "361×5=1805"
Replace-Text "889×5=4445" "161×9=1449"
Replace-Text "893×5=4465" "666×4=2664"
Replace-Text "684×3=2052" "864×4=3456"
Replace-Text "952×7=6664" "795×6=4770"
Replace-Text "242×3=726" "558×8=4464"
Replace-Text "419×8=3352" "873×9=7857"
Replace-Text "949×4=3796" "270×2=540"
Replace-Text "178×6=1068" "445×6=2670"
Replace-Text "522×5=2610" "432×8=3456"
Replace-Text "262×4=1048" "916×8=7328"

Write-Output "Done"
